$d = $word.ActiveDocument

$d.Content.Find.Execute("293÷8=36, 5", $true, $false, $false, $false, $false, $true, 1, $false, "123÷6=20, 3", 2) | Out-Null
$d.Content.Find.Execute("133÷2=66, 1", $true, $false, $false, $false, $false, $true, 1, $false, "995÷3=331, 2", 2) | Out-Null
$d.Content.Find.Execute("566÷9=62, 8", $true, $false, $false, $false, $false, $true, 1, $false, "320÷3=106, 2", 2) | Out-Null
$d.Content.Find.Execute("811÷3=270, 1", $true, $false, $false, $false, $false, $true, 1, $false, "131÷7=18, 5", 2) | Out-Null
$d.Content.Find.Execute("295÷7=42, 1", $true, $false, $false, $false, $false, $true, 1, $false, "570÷8=71, 2", 2) | Out-Null
$d.Content.Find.Execute("361÷6=60, 1", $true, $false, $false, $false, $false, $true, 1, $false, "356÷9=39, 5", 2) | Out-Null
$d.Content.Find.Execute("102÷9=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "385÷4=96, 1", 2) | Out-Null
$d.Content.Find.Execute("269÷4=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "758÷6=126, 2", 2) | Out-Null
$d.Content.Find.Execute("266÷6=44, 2", $true, $false, $false, $false, $false, $true, 1, $false, "623÷6=103, 5", 2) | Out-Null
$d.Content.Find.Execute("773÷9=85, 8", $true, $false, $false, $false, $false, $true, 1, $false, "743÷2=371, 1", 2) | Out-Null
$d.Content.Find.Execute("626÷2=313, 0", $true, $false, $false, $false, $false, $true, 1, $false, "555÷5=111, 0", 2) | Out-Null
$d.Content.Find.Execute("355÷3=118, 1", $true, $false, $false, $false, $false, $true, 1, $false, "513÷3=171, 0", 2) | Out-Null
$d.Content.Find.Execute("138÷9=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "305÷4=76, 1", 2) | Out-Null
$d.Content.Find.Execute("112÷5=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "495÷8=61, 7", 2) | Out-Null
$d.Content.Find.Execute("985÷3=328, 1", $true, $false, $false, $false, $false, $true, 1, $false, "367÷4=91, 3", 2) | Out-Null
$d.Content.Find.Execute("802÷7=114, 4", $true, $false, $false, $false, $false, $true, 1, $false, "477÷9=53, 0", 2) | Out-Null
$d.Content.Find.Execute("661÷2=330, 1", $true, $false, $false, $false, $false, $true, 1, $false, "578÷3=192, 2", 2) | Out-Null
$d.Content.Find.Execute("548÷3=182, 2", $true, $false, $false, $false, $false, $true, 1, $false, "928÷4=232, 0", 2) | Out-Null
$d.Content.Find.Execute("446÷8=55, 6", $true, $false, $false, $false, $false, $true, 1, $false, "551÷8=68, 7", 2) | Out-Null
$d.Content.Find.Execute("659÷6=109, 5", $true, $false, $false, $false, $false, $true, 1, $false, "382÷5=76, 2", 2) | Out-Null
$d.Content.Find.Execute("369÷6=61, 3", $true, $false, $false, $false, $false, $true, 1, $false, "112÷3=37, 1", 2) | Out-Null
$d.Content.Find.Execute("849÷4=212, 1", $true, $false, $false, $false, $false, $true, 1, $false, "834÷6=139, 0", 2) | Out-Null
$d.Content.Find.Execute("108÷3=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "608÷3=202, 2", 2) | Out-Null
$d.Content.Find.Execute("687÷8=85, 7", $true, $false, $false, $false, $false, $true, 1, $false, "604÷7=86, 2", 2) | Out-Null
$d.Content.Find.Execute("418÷8=52, 2", $true, $false, $false, $false, $false, $true, 1, $false, "689÷4=172, 1", 2) | Out-Null
